$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVIV_CheckoutPage")

# D6 ("Address" test data) and D8 ("PhNo" test data) hold digit-only strings.
# Prefix with an apostrophe (exactly like typing '12345 into Excel) so the
# engine stores them as text (shared string) instead of auto-converting to
# a number, then reset the style back to Normal so no quote-prefix style
# lingers on the cell (matches the original s="0" plain-text shared string).
$ws.Range("D6").Value = "'5524032318"
$ws.Range("D6").Style = "Normal"

$ws.Range("D8").Value = "'4199964175"
$ws.Range("D8").Style = "Normal"

$ws.Range("D13").Value = "ORDER NUMBER: 1037"
